# Added October NFTF Link
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 15: month label in column A, newsletter hyperlink in column B
$ws.Range("A15").Value = "October 2019"
$ws.Range("B15").Value = "https://myemail.constantcontact.com/News-From-The-Forest---October.html?soid=1102494320279&aid=t-ew4tkMBqU"

# Match the formatting used by the rows above it
$ws.Range("A15").NumberFormat = $ws.Range("A14").NumberFormat

$ws.Hyperlinks.Add($ws.Range("B15"), "https://myemail.constantcontact.com/News-From-The-Forest---October.html?soid=1102494320279&aid=t-ew4tkMBqU") | Out-Null
$ws.Range("B15").Style = "Hyperlink"

# Matches the selection left behind in the authored workbook
$ws.Range("B21").Select() | Out-Null
